# "Actualización automática desde WSL"
# The sensor-logger workbook is re-synced: the last week's worth of
# timestamps in column A (rows 14-20) get refreshed to the latest
# reading times, and the active selection moves to the first empty
# row below the data (A21) instead of the mid-sheet cell left
# selected previously.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the timestamp (serial date/time) values for rows 14-20.
$ws.Range("A14").Value = 45875.70853009259
$ws.Range("A15").Value = 45875.750196759262
$ws.Range("A16").Value = 45875.791863425926
$ws.Range("A17").Value = 45875.83353009259
$ws.Range("A18").Value = 45875.875196759262
$ws.Range("A19").Value = 45875.916863425926
$ws.Range("A20").Value = 45875.95853009259

# Move the active selection to A21 (first row after the data),
# matching where the app would place the cursor after appending.
$ws.Range("A21").Select()
